$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.644.27'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').Value = '2.292.75'
$ws.Range('E3').Value = '  -0.77%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = "'96.45"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.28%  '

$ws.Range('D6').Value = "'268.80"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.07%  '

$ws.Range('D7').Value = "'0.623"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.65%  '

$ws.Range('D8').Value = "'0.999"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').Value = "'0.611"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.34%  '

$ws.Range('D10').Value = "'45.51"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.23%  '

$ws.Range('D11').Value = "'0.0936"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.22%  '

$ws.Range('D12').Value = "'7.94"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.33%  '

$ws.Range('E13').Value = '  +2.14%  '

$ws.Range('D14').Value = '2.634.38'
$ws.Range('E14').Value = '  -0.82%  '

$ws.Range('D15').Value = "'15.40"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.66%  '

$ws.Range('D16').Value = "'0.849"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.59%  '

$ws.Range('D17').Value = '2.288.13'
$ws.Range('E17').Value = '  -1.69%  '

$ws.Range('D18').Value = '43.624.57'
$ws.Range('E18').Value = '  -0.09%  '

$ws.Range('D19').Value = "'0.0000111"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.50%  '

$ws.Range('D20').Value = "'6.20"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.54%  '

$ws.Range('D21').Value = "'72.13"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.04%  '

$ws.Range('D22').Value = "'2.55"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +12.37%  '

$ws.Range('D23').Value = "'232.84"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.08%  '

$ws.Range('D24').Value = "'9.12"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.29%  '

$ws.Range('D25').Value = "'2.64"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.24%  '

$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('D27').Value = "'11.26"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '

$ws.Range('E28').Value = '  +2.60%  '

$ws.Range('D29').Value = "'39.85"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.40%  '

$ws.Range('D30').Value = "'2.28"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.15%  '

$ws.Range('D31').Value = "'175.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.03%  '

$ws.Range('D32').Value = "'21.86"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.12%  '

$ws.Range('D33').Value = "'0.0898"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.81%  '

$ws.Range('D34').Value = "'5.40"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.40%  '

$ws.Range('E35').Value = '  -0.29%  '

$ws.Range('E36').Value = '  -1.30%  '

$ws.Range('D37').Value = "'0.0353"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.85%  '

$ws.Range('D38').Value = "'4.38"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.26%  '

$ws.Range('D39').Value = "'3.39"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.00%  '

$ws.Range('D40').Value = "'0.240"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.72%  '

$ws.Range('D41').Value = "'2.34"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.54%  '

$ws.Range('D42').Value = "'12.34"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.90%  '

$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'1.35"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.09%  '

$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').Value = "'65.31"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.87%  '

$ws.Range('E45').Value = '  -2.57%  '

$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = "'5.16"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.16%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = "'0.102"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '

$ws.Range('D48').Value = "'97.46"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.63%  '

$ws.Range('E49').Value = '  -0.73%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = "'1.52"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.91%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.515.75'
$ws.Range('E51').Value = '  -0.72%  '
